$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 565.25
$ws.Range("I111").Value = 517.5714
$ws.Range("J111").Value = 899
$ws.Range("K111").Value = 1552.7142
$ws.Range("L111").Value = 2697
$ws.Range("M111").Value = 1514.2858
$ws.Range("N111").Value = -8831
$ws.Range("H113").Value = 6181.273
$ws.Range("I113").Value = 4600
$ws.Range("J113").Value = 7499
$ws.Range("K113").Value = 4600
$ws.Range("L113").Value = 7499
$ws.Range("M113").Value = -1346
$ws.Range("N113").Value = -14007
$ws.Range("H137").Value = 2696.8572
$ws.Range("I137").Value = 2696.8572
$ws.Range("K137").Value = 8090.571599999999
$ws.Range("M137").Value = -5540.571599999999
$ws.Range("H138").Value = 2546.7693
$ws.Range("I138").Value = 656.4545000000001
$ws.Range("K138").Value = 1969.3635
$ws.Range("M138").Value = 3170.6365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 291.64285
$ws.Range("I4").Value = 239.75
$ws.Range("J4").Value = 360.83334
$ws.Range("K4").Value = 239.75
$ws.Range("L4").Value = 360.83334
$ws.Range("M4").Value = -123.75
$ws.Range("N4").Value = -592.83334
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H23").Value = 12000
$ws.Range("J23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("N23").Value = -12518
$ws.Range("H36").Value = 8026
$ws.Range("I36").Value = 8026
$ws.Range("K36").Value = 8026
$ws.Range("M36").Value = -7680
$ws.Range("H37").Value = 9150
$ws.Range("I37").Value = 8300
$ws.Range("K37").Value = 8300
$ws.Range("M37").Value = -8027
$ws.Range("H55").Value = 19997.5
$ws.Range("J55").Value = 19997.5
$ws.Range("L55").Value = 19997.5
$ws.Range("N55").Value = -20627.5
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H95").Value = 8440.6
$ws.Range("J95").Value = 8440.6
$ws.Range("L95").Value = 8440.6
$ws.Range("N95").Value = -13932.6
$ws.Range("H132").Value = 4317
$ws.Range("I132").Value = 4317
$ws.Range("K132").Value = 12951
$ws.Range("M132").Value = -10421
$ws.Range("H140").Value = 365000
$ws.Range("J140").Value = 365000
$ws.Range("L140").Value = 365000
$ws.Range("N140").Value = -375360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H15").Value = 23986
$ws.Range("J15").Value = 25482.5
$ws.Range("L15").Value = 25482.5
$ws.Range("N15").Value = -25936.5
$ws.Range("H19").Value = 32487.5
$ws.Range("I19").Value = 29975
$ws.Range("K19").Value = 29975
$ws.Range("M19").Value = -29802
$ws.Range("H35").Value = 14998.333
$ws.Range("J35").Value = 14998.333
$ws.Range("L35").Value = 14998.333
$ws.Range("N35").Value = -15618.333
$ws.Range("H82").Value = 22765.375
$ws.Range("J82").Value = 29998.75
$ws.Range("L82").Value = 29998.75
$ws.Range("N82").Value = -30764.75
$ws.Range("H85").Value = 22765.375
$ws.Range("J85").Value = 29998.75
$ws.Range("L85").Value = 29998.75
$ws.Range("N85").Value = -32650.75
$ws.Range("H94").Value = 10010
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 10010
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 10010
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -10912
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H137").Value = 41666
$ws.Range("I137").Value = 25000
$ws.Range("K137").Value = 25000
$ws.Range("M137").Value = -19900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("H62").Value = 5250
$ws.Range("J62").Value = 5250
$ws.Range("L62").Value = 5250
$ws.Range("N62").Value = -6498
$ws.Range("H65").Value = 5250
$ws.Range("J65").Value = 5250
$ws.Range("L65").Value = 26250
$ws.Range("N65").Value = -32490
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
$ws.Range("H99").Value = 1804.6666
$ws.Range("I99").Value = 1804.6666
$ws.Range("K99").Value = 1804.6666
$ws.Range("M99").Value = -306.6666
$ws.Range("H126").Value = 1804.6666
$ws.Range("I126").Value = 1804.6666
$ws.Range("K126").Value = 5413.9998
$ws.Range("M126").Value = -2943.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 502.36365
$ws.Range("J2").Value = 371.66666
$ws.Range("L2").Value = 371.66666
$ws.Range("N2").Value = -597.66666
$ws.Range("H18").Value = 39000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 39000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 39000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -39586
$ws.Range("H27").Value = 5625
$ws.Range("J27").Value = 5625
$ws.Range("L27").Value = 5625
$ws.Range("N27").Value = -5957
$ws.Range("H43").Value = 3436.8572
$ws.Range("I43").Value = 2343
$ws.Range("K43").Value = 2343
$ws.Range("M43").Value = -2192
$ws.Range("H57").Value = 16248.75
$ws.Range("J57").Value = 19998.334
$ws.Range("L57").Value = 19998.334
$ws.Range("N57").Value = -21638.334
$ws.Range("H80").Value = 1518.3334
$ws.Range("I80").Value = 1518.3334
$ws.Range("K80").Value = 1518.3334
$ws.Range("M80").Value = -520.3334
$ws.Range("H83").Value = 1518.3334
$ws.Range("I83").Value = 1518.3334
$ws.Range("K83").Value = 7591.666999999999
$ws.Range("M83").Value = -2599.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15500.125
$ws.Range("I7").Value = 15500.125
$ws.Range("K7").Value = 15500.125
$ws.Range("M7").Value = -15388.125
$ws.Range("H46").Value = 4553.8
$ws.Range("I46").Value = 8391.6
$ws.Range("J46").Value = 2634.9
$ws.Range("K46").Value = 8391.6
$ws.Range("L46").Value = 2634.9
$ws.Range("M46").Value = -8203.6
$ws.Range("N46").Value = -3010.9
$ws.Range("H126").Value = 15500.125
$ws.Range("I126").Value = 15500.125
$ws.Range("K126").Value = 46500.375
$ws.Range("M126").Value = -44030.375
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 39274.43
$ws.Range("I126").Value = 34603.31
$ws.Range("K126").Value = 103809.93
$ws.Range("M126").Value = -101339.93
